$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CardInfo")

# Set custom column widths for E, F, G, H (values chosen so the stored/quantized
# OOXML width lands on the closest achievable value to the target widths
# 20.5, 19.875, 55, 75.875).
$ws.Columns.Item(5).ColumnWidth = 19.8
$ws.Columns.Item(6).ColumnWidth = 19.15
$ws.Columns.Item(7).ColumnWidth = 54.3
$ws.Columns.Item(8).ColumnWidth = 75.15

# Populate the new "CardContent" (English) column G for every data row,
# and fix the CardContentKor text for row 18 (Shadow Reaper).
$ws.Range('G2').Value = 'Deal 5 physical damage.'
$ws.Range('G3').Value = 'Gain 5 extra armor.'
$ws.Range('G4').Value = 'Heal 10 hp. Remove bleeding.'
$ws.Range('G5').Value = 'Deal 5 magical damage.'
$ws.Range('G6').Value = 'Draw 2 cards.'
$ws.Range('G7').Value = 'Draw 1 card. Remove Exhausted.'
$ws.Range('G8').Value = 'Deal 7 physical damage with 2-turn bleeding.'
$ws.Range('G9').Value = 'Deal 4 physical damage with 2-turn bleeding.'
$ws.Range('G10').Value = 'Deal 5 true damage.'
$ws.Range('G11').Value = 'Lose 10 hp. Gain 15 extra armor.'
$ws.Range('G12').Value = 'Deal 7 magical damage with 2-turn posion.'
$ws.Range('G13').Value = 'Deal 3 magical damage to whole enemy with 2-turn dizziness.'
$ws.Range('G14').Value = 'Find and bring 1 card from used card. The card costs 0.'
$ws.Range('G15').Value = 'Deal 15 magical damage. If any enemy killed, take Stealth stance for 1 turn.'
$ws.Range('G16').Value = 'Gain 1 extra composure.'
$ws.Range('G17').Value = 'Gain extra 25 golds and heal 10 hp for every kill.'
$ws.Range('G18').Value = 'Deal 4 magical damage to random target every end of player''s turn.'
$ws.Range('G19').Value = 'Deal 3 physical damage.'
$ws.Range('G20').Value = 'Deal 4 physical damage. Draw 1 card.'
$ws.Range('G21').Value = 'Deal 6 magical damage to whole enemy.'
$ws.Range('G22').Value = 'Deal 4 damage to random target for 8 times.'
$ws.Range('G23').Value = 'Fall into Exhausted and Dizziness in this turn. Gain 1 extra intelligence for next turn.'
$ws.Range('G24').Value = 'Give every enemy 2-turn Exhausted and Paranoia.'
$ws.Range('G25').Value = 'Gain 10 extra armor. Take CounterAttack stance for 1 turn.'
$ws.Range('G26').Value = 'Fill your hand full of "Single Shot" cards. Fall into Exhausted and Dizziness in next turn.'
$ws.Range('G27').Value = 'Deal 1 magical damage to whole enemy after using every action card.'
$ws.Range('G28').Value = 'Gain 5 extra strength. Fall into Paranoia.'
$ws.Range('G29').Value = 'Deal extra 1 magical damage for every physical damage.'

# Row 18 (Shadow Reaper) Korean description was also changed.
$ws.Range('H18').Value = '턴 종료 시 무작위 적에게 마법피해4'

# Update the active selection to match the saved view state.
$ws.Activate()
$ws.Range("G15").Select() | Out-Null
